# Adds a new "DGC" worksheet at the end of the workbook with a small
# 3-column table (COMARCA / TEMÁTICA / PROBLEMA) describing recurring
# DGC issues for Dom Pedrito.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet as the LAST tab -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "DGC"

# --- 2. Table contents -------------------------------------------------
$ws.Range("A1").Value = "COMARCA"
$ws.Range("B1").Value = "TEMÁTICA"
$ws.Range("C1").Value = "PROBLEMA"

$ws.Range("A2").Value = "Dom Pedrito"
$ws.Range("B2").Value = "MOT-Limpeza"
$ws.Range("C2").Value = "Demora nos Atestes MOT"

$ws.Range("A3").Value = "Dom Pedrito"
$ws.Range("B3").Value = "MOT-Limpeza"
$ws.Range("C3").Value = "Demora nos Atestes MOT/não responde TEAMS"

$ws.Range("A4").Value = "Dom Pedrito"
$ws.Range("B4").Value = "MOT-VIG"
$ws.Range("C4").Value = "Demora nos Atestes MOT"

$ws.Range("A5").Value = "Dom Pedrito"
$ws.Range("B5").Value = "MAOBRAS/MATIC"
$ws.Range("C5").Value = "Atraso no ateste."

# --- 3. Column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.42578125
$ws.Columns.Item(2).ColumnWidth = 15.85546875
$ws.Columns.Item(3).ColumnWidth = 88.140625

# --- 4. Header row formatting (bold, gray fill, centered, bordered) -------
$header = $ws.Range("A1:C1")
$header.Font.Name = "Aptos Narrow"
$header.Font.Size = 11
$header.Font.Bold = $true
$header.Font.Color = 0
$header.Interior.Pattern = 1
$header.Interior.Color = 14277081
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# --- 5. Data rows formatting (plain font, bordered) ------------------------
$body = $ws.Range("A2:C5")
$body.Font.Name = "Aptos Narrow"
$body.Font.Size = 11
$body.Font.Bold = $false
$body.Font.Color = 0
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2

Write-Host "DGC sheet added"
